$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.931.60'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.767.67'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.03'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4732'
$ws.Range('E7').Value = '  +3.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3527'
$ws.Range('E8').Value = '  -1.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.62'
$ws.Range('E9').Value = '  +4.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07386'
$ws.Range('E10').Value = '  -1.40%  '
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('E13').Value = '  -0.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.017'
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.188'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').Value = '1.769.40'
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.18'
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06416'
$ws.Range('E19').Value = '  -0.32%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.794'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = '27.977.87'
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('E24').Value = '  -1.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.159'
$ws.Range('E25').Value = '  +3.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.05'
$ws.Range('E26').Value = '  +0.76%  '
$ws.Range('E27').Value = '  -0.99%  '
$ws.Range('D28').Value = '1.963.50'
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.207'
$ws.Range('E29').Value = '  +1.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.35'
$ws.Range('E30').Value = '  -1.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.074'
$ws.Range('E31').Value = '  -2.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09380'
$ws.Range('E32').Value = '  +1.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.658'
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.69'
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06113'
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02262'
$ws.Range('E37').Value = '  -1.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2072'
$ws.Range('E38').Value = '  -0.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.905'
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6175'
$ws.Range('E40').Value = '  -2.35%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.192'
$ws.Range('E41').Value = '  +0.44%  '
$ws.Range('B42').Value = 'WEMIXTOKEN'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.447'
$ws.Range('E42').Value = '  +4.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.756'
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.18'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5803'
$ws.Range('E46').Value = '  -2.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '123.96'
$ws.Range('E47').Value = '  +0.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.933'
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.128'
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06810'
$ws.Range('E50').Value = '  -1.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.06'
$ws.Range('E51').Value = '  -0.39%  '
